# Lab5 "admin" sheet register: fix row 17 (B3 correction per commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B17 previously carried a manually-applied fill style; remove that
# direct formatting so the cell reverts to the default/normal style.
$ws.Range("B17").ClearFormats()

# C17's numeric value was corrected from 123456 to 12345.
$ws.Range("C17").Value = 12345

# Leave the active selection on C17, matching the last-edited cell.
$ws.Range("C17").Select()
